# Add a new test-data row (row 5) to the recruiter-tracking sheet, add the
# matching mailto: hyperlink for the recruiter e-mail cell, resize the last
# two columns to fit the new (longer) content, and leave the selection where
# the author last clicked.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of test data ------------------------------------------------
$ws.Range("A5").Value = "Company D"
$ws.Range("B5").Value = "AI researcher"
$ws.Range("C5").Value = "unknown@gmail.com"
$ws.Range("D5").Value = "non_existant_template.pdf"

# --- Hyperlink the new recruiter e-mail, matching the style used by the
#     other rows in column C (the built-in "Hyperlink" cell style) --------
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:unknown@gmail.com")
$ws.Range("C5").Style = "Hyperlink"

# --- Widen columns C & D so the new, longer values fit (mirrors an
#     Excel "AutoFit" on those columns) -----------------------------------
$ws.Columns("C").ColumnWidth = 18.053385416666668
$ws.Columns("D").ColumnWidth = 22.608072916666668

# --- Leave the selection where the author ended up ------------------------
$ws.Range("C7").Select()
